# Apply the update described by the diff:
#  - Row 2: change O2, P2 odds
#  - Row 3: change G3, H3, I3, J3, L3, M3, N3, U3, V3 odds
#  - Add two new data rows (5 and 6) for Venezuela - Liga Futve fixtures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Cells.Item(2, 15).Value = 1.36   # O2
$ws.Cells.Item(2, 16).Value = 3      # P2

# --- Row 3 updates ---
$ws.Cells.Item(3, 7).Value  = 8.5    # G3
$ws.Cells.Item(3, 8).Value  = 5      # H3
$ws.Cells.Item(3, 9).Value  = 1.36   # I3
$ws.Cells.Item(3, 10).Value = 8.5    # J3
$ws.Cells.Item(3, 12).Value = 1.83   # L3
$ws.Cells.Item(3, 13).Value = 1.06   # M3
$ws.Cells.Item(3, 14).Value = 10     # N3
$ws.Cells.Item(3, 21).Value = 2.25   # U3
$ws.Cells.Item(3, 22).Value = 1.57   # V3

# --- New rows 5 and 6 ---
# Columns A..BD, in order.
$row5Values = @("Qa7iAtsI", "13/11/2024", "19:30", "VENEZUELA - LIGA FUTVE", "Carabobo", "Monagas", 1.62, 3.5, 5.3, 2.2, 2.07, 5.5, 1.05, 6.2, 1.37, 2.62, 2.07, 1.6, 1.44, 2.42, 2.05, 1.6, 5.4, 6.6, 8.5, 11.5, 15, 37, 7.9, 6.9, 20, 120, 11.75, 30, 18, 110, 65, 75, 3.3, 7.9, 20, 27, 70, 2.4, 8.25, 90, 6.7, 32, 40, 200, 300, 500, 900, 350, 51, 51)
$row6Values = @("MR008KBU", "13/11/2024", "19:30", "VENEZUELA - LIGA FUTVE", "La Guaira", "Estudiantes Merida", 2.02, 3.4, 3.3, 2.62, 2.15, 3.65, 1.01, 8.1, 1.22, 3.4, 1.7, 1.93, 1.33, 3.04, 1.57, 2.1, 8.5, 10.5, 8.5, 19, 15.5, 23, 12, 6.8, 12.5, 45, 12, 19.5, 11.25, 45, 26, 29, 4.05, 10.5, 17.5, 40, 65, 2.87, 6.6, 50, 5.3, 17, 22, 80, 100, 250, 300, 200, 51, 51)

for ($i = 0; $i -lt $row5Values.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5Values[$i]
}

for ($i = 0; $i -lt $row6Values.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6Values[$i]
}
